$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Stock"

$ws.Columns.Item(1).ColumnWidth = 2.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 5.5
$ws.Columns.Item(3).ColumnWidth = 4.666666666666667
$ws.Columns.Item(4).ColumnWidth = 5.0
